$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "duration" header to "duration (hrs)" and fill in the new
# duration values for every activity row.
$ws.Range("B3").Value = "duration (hrs)"
$ws.Range("B4").Value = 0.25
$ws.Range("B5").Value = 1
$ws.Range("B6").Value = 2
$ws.Range("B7").Value = 0.5
$ws.Range("B8").Value = 0.5
$ws.Range("B9").Value = 0.5
$ws.Range("B10").Value = 0.5
$ws.Range("B11").Value = 1
$ws.Range("B12").Value = 1.5
$ws.Range("B13").Value = 1
$ws.Range("B14").Value = 1.5
$ws.Range("B15").Value = 0.25
$ws.Range("B16").Value = 0.25
$ws.Range("B17").Value = 0.5
$ws.Range("B18").Value = 0.25
$ws.Range("B19").Value = 0.5
$ws.Range("B20").Value = 0.5
$ws.Range("B21").Value = 1
$ws.Range("B22").Value = 0.5
$ws.Range("B23").Value = 1

# Widen the new duration column like the rest of the table.
$ws.Columns.Item(2).ColumnWidth = 12.85546875

# Add the "Crib Assembly" title above the table, spanning C1:D1, styled as
# a bold, red, 18pt heading.
$ws.Range("C1").Value = "Crib Assembly"
$ws.Range("C1:D1").Merge()
$ws.Range("C1:D1").Font.Bold = $true
$ws.Range("C1:D1").Font.Size = 18
$ws.Range("C1:D1").Font.Color = 255
$ws.Rows.Item(1).RowHeight = 23.25

# Set the page to print in portrait orientation.
$ws.PageSetup.Orientation = 1

# Match the saved selection from the edited workbook.
$ws.Range("F11").Select()
